$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to remain
# text (matching the original inlineStr/shared-string cell type) by
# temporarily applying a text number format, then resetting the style so
# no extra formatting is left behind.
$textForceRefs = @(
    "D4", "D6", "D14", "D15", "D20", "D21", "D22", "D25", "D26", "D28", "D29", "D30", "D32", "D34", "D37", "D38", "D39", "D40", "D42", "D44", "D45", "D46", "D48", "D51"
)
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply all cell value updates per the diff
$ws.Range("D2").Value = "65.111.11"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "3.528.78"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "134.72"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D7").Value = "3.524.30"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "4.130.07"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "27.75"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "0.0000183"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "3.526.59"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "65.102.04"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "14.47"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "5.73"
$ws.Range("D22").Value = "392.94"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "3.671.42"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").Value = "74.76"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -3.83%  "
$ws.Range("D28").Value = "7.76"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "1.57"
$ws.Range("E29").Value = "  +9.30%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "8.41"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "3.532.49"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").Value = "24.22"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").Value = "5.30"
$ws.Range("E37").Value = "  +5.35%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.58"
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "7.00"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "168.68"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "0.824"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  +3.97%  "
$ws.Range("D44").Value = "25.87"
$ws.Range("E44").Value = "  -4.13%  "
$ws.Range("D45").Value = "42.94"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").Value = "1.66"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "2.410.91"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("D51").Value = "0.907"
$ws.Range("E51").Value = "  +5.61%  "

# Reset style on the text-forced cells so no residual number-format
# styling remains on the cell (keeps cell style identical to original).
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).Style = "Normal"
}
